$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K10").Value = 5.479641
$ws.Range("L10").Value = 5.479641
$ws.Range("M10").Value = 5.479641
$ws.Range("N10").Value = 8.225541
$ws.Range("O10").Value = 8.225541
$ws.Range("P10").Value = 9.753021
$ws.Range("Q10").Value = 9.753021
$ws.Range("R10").Value = 11.258711
$ws.Range("S10").Value = 11.258711
